$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6 (anchor G6=4564) on ALC
$ws.Range("H6").Value2 = 354.625
$ws.Range("I6").Value2 = 354.625
$ws.Range("K6").Value2 = 1063.875
$ws.Range("M6").Value2 = -951.875

# Row 8 (anchor G8=4565) on ALC
$ws.Range("H8").Value2 = 207.33333
$ws.Range("I8").Value2 = 48.8
$ws.Range("J8").Value2 = 1000
$ws.Range("K8").Value2 = 146.4
$ws.Range("L8").Value2 = 3000
$ws.Range("M8").Value2 = -7.399999999999977
$ws.Range("N8").Value2 = -3278

# Row 31 (anchor G31=4576) on ALC
$ws.Range("H31").Value2 = 228.66667
$ws.Range("I31").Value2 = 228.66667
$ws.Range("K31").Value2 = 686.00001
$ws.Range("M31").Value2 = -456.00001

# Row 39 (anchor G39=4603) on ALC
$ws.Range("H39").Value2 = 625.2727
$ws.Range("I39").Value2 = 504.66666
$ws.Range("J39").Value2 = 770
$ws.Range("K39").Value2 = 1513.99998
$ws.Range("L39").Value2 = 2310
$ws.Range("M39").Value2 = -1217.99998
$ws.Range("N39").Value2 = -2902

# Row 48 (anchor G48=4587) on ALC
$ws.Range("H48").Value2 = 1019
$ws.Range("I48").Value2 = 0
$ws.Range("J48").Value2 = 1019
$ws.Range("K48").Value2 = 0
$ws.Range("L48").Value2 = 3057
$ws.Range("M48").Value2 = ""
$ws.Range("N48").Value2 = -3641

# Row 55 (anchor G55=5517) on ALC
$ws.Range("H55").Value2 = 189.64285
$ws.Range("J55").Value2 = 234.14285
$ws.Range("L55").Value2 = 234.14285
$ws.Range("N55").Value2 = -662.14285

# Row 56 (anchor G56=4587) on ALC
$ws.Range("H56").Value2 = 1019
$ws.Range("I56").Value2 = 0
$ws.Range("J56").Value2 = 1019
$ws.Range("K56").Value2 = 0
$ws.Range("L56").Value2 = 3057
$ws.Range("M56").Value2 = ""
$ws.Range("N56").Value2 = -4125

# Row 88 (anchor G88=12608) on ALC
$ws.Range("H88").Value2 = 3452.9092
$ws.Range("J88").Value2 = 3415.2222
$ws.Range("L88").Value2 = 3415.2222
$ws.Range("N88").Value2 = -4227.2222

# Row 91 (anchor G91=12608) on ALC
$ws.Range("H91").Value2 = 3452.9092
$ws.Range("J91").Value2 = 3415.2222
$ws.Range("L91").Value2 = 3415.2222
$ws.Range("N91").Value2 = -6223.2222

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (anchor G32=44147) on ARM
$ws.Range("H32").Value2 = 4335.952
$ws.Range("I32").Value2 = 3404.8333
$ws.Range("K32").Value2 = 3404.8333
$ws.Range("M32").Value2 = -3117.8333

# Row 92 (anchor G92=18050) on ARM
$ws.Range("H92").Value2 = 0
$ws.Range("J92").Value2 = 0
$ws.Range("L92").Value2 = 0
$ws.Range("N92").Value2 = ""

# Row 122 (anchor G122=36168) on ARM
$ws.Range("H122").Value2 = 1669
$ws.Range("I122").Value2 = 1669
$ws.Range("K122").Value2 = 5007
$ws.Range("M122").Value2 = -2557

# Row 132 (anchor G132=43997) on ARM
$ws.Range("H132").Value2 = 1574.75
$ws.Range("I132").Value2 = 1733.3334
$ws.Range("J132").Value2 = 1099
$ws.Range("K132").Value2 = 5200.0002
$ws.Range("L132").Value2 = 3297
$ws.Range("M132").Value2 = -2670.0002
$ws.Range("N132").Value2 = -8357

$ws = $wb.Worksheets.Item("BSM")
# Row 22 (anchor G22=5092) on BSM
$ws.Range("H22").Value2 = 523.0833
$ws.Range("I22").Value2 = 523.0833
$ws.Range("J22").Value2 = 0
$ws.Range("K22").Value2 = 523.0833
$ws.Range("L22").Value2 = 0
$ws.Range("M22").Value2 = -350.0833
$ws.Range("N22").Value2 = ""

# Row 92 (anchor G92=18033) on BSM
$ws.Range("H92").Value2 = 39999
$ws.Range("J92").Value2 = 39999
$ws.Range("L92").Value2 = 39999
$ws.Range("N92").Value2 = -44991

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (anchor G16=27691) on CRP
$ws.Range("H16").Value2 = 1103.6666
$ws.Range("I16").Value2 = 1103.6666
$ws.Range("J16").Value2 = 0
$ws.Range("K16").Value2 = 1103.6666
$ws.Range("L16").Value2 = 0
$ws.Range("M16").Value2 = -816.6666
$ws.Range("N16").Value2 = ""

# Row 63 (anchor G63=10604) on CRP
$ws.Range("H63").Value2 = 0
$ws.Range("J63").Value2 = 0
$ws.Range("L63").Value2 = 0
$ws.Range("N63").Value2 = ""

# Row 66 (anchor G66=10604) on CRP
$ws.Range("H66").Value2 = 0
$ws.Range("J66").Value2 = 0
$ws.Range("L66").Value2 = 0
$ws.Range("N66").Value2 = ""

# Row 113 (anchor G113=27691) on CRP
$ws.Range("H113").Value2 = 1103.6666
$ws.Range("I113").Value2 = 1103.6666
$ws.Range("J113").Value2 = 0
$ws.Range("K113").Value2 = 1103.6666
$ws.Range("L113").Value2 = 0
$ws.Range("M113").Value2 = 1066.3334
$ws.Range("N113").Value2 = ""

$ws = $wb.Worksheets.Item("CUL")
# Row 7 (anchor G7=4728) on CUL
$ws.Range("H7").Value2 = 227
$ws.Range("I7").Value2 = 90.5
$ws.Range("J7").Value2 = 336.2
$ws.Range("K7").Value2 = 271.5
$ws.Range("L7").Value2 = 1008.6
$ws.Range("M7").Value2 = -159.5
$ws.Range("N7").Value2 = -1232.6

# Row 11 (anchor G11=4745) on CUL
$ws.Range("H11").Value2 = 2066.75
$ws.Range("I11").Value2 = 2220
$ws.Range("J11").Value2 = 1300.5
$ws.Range("K11").Value2 = 6660
$ws.Range("L11").Value2 = 3901.5
$ws.Range("M11").Value2 = -6520
$ws.Range("N11").Value2 = -4181.5

# Row 14 (anchor G14=12886) on CUL
$ws.Range("H14").Value2 = 2524.25
$ws.Range("I14").Value2 = 2524.25
$ws.Range("K14").Value2 = 7572.75
$ws.Range("M14").Value2 = -7399.75

# Row 38 (anchor G38=4860) on CUL
$ws.Range("H38").Value2 = 108.15385
$ws.Range("I38").Value2 = 41.285713
$ws.Range("J38").Value2 = 186.16667
$ws.Range("K38").Value2 = 123.857139
$ws.Range("L38").Value2 = 558.50001
$ws.Range("M38").Value2 = 223.142861
$ws.Range("N38").Value2 = -1252.50001

# Row 41 (anchor G41=4700) on CUL
$ws.Range("H41").Value2 = 179.5
$ws.Range("I41").Value2 = 179.5
$ws.Range("J41").Value2 = 0
$ws.Range("K41").Value2 = 538.5
$ws.Range("L41").Value2 = 0
$ws.Range("M41").Value2 = -200.5
$ws.Range("N41").Value2 = ""

# Row 98 (anchor G98=19843) on CUL
$ws.Range("H98").Value2 = 300
$ws.Range("J98").Value2 = 300
$ws.Range("L98").Value2 = 900
$ws.Range("N98").Value2 = -3896

# Row 103 (anchor G103=19839) on CUL
$ws.Range("H103").Value2 = 553.8182
$ws.Range("I103").Value2 = 109.14286
$ws.Range("K103").Value2 = 327.42858
$ws.Range("M103").Value2 = 551.57142

# Row 108 (anchor G108=27853) on CUL
$ws.Range("H108").Value2 = 323
$ws.Range("I108").Value2 = 323
$ws.Range("K108").Value2 = 969
$ws.Range("M108").Value2 = 1911

$ws = $wb.Worksheets.Item("GSM")
# Row 80 (anchor G80=12521) on GSM
$ws.Range("H80").Value2 = 2750
$ws.Range("I80").Value2 = 2666.6667
$ws.Range("J80").Value2 = 2833.3333
$ws.Range("K80").Value2 = 2666.6667
$ws.Range("L80").Value2 = 2833.3333
$ws.Range("M80").Value2 = -1668.6667
$ws.Range("N80").Value2 = -4829.3333

# Row 83 (anchor G83=12521) on GSM
$ws.Range("H83").Value2 = 2750
$ws.Range("I83").Value2 = 2666.6667
$ws.Range("J83").Value2 = 2833.3333
$ws.Range("K83").Value2 = 13333.3335
$ws.Range("L83").Value2 = 14166.6665
$ws.Range("M83").Value2 = -8341.333500000001
$ws.Range("N83").Value2 = -24150.6665

# Row 132 (anchor G132=44008) on GSM
$ws.Range("H132").Value2 = 4999.3335
$ws.Range("I132").Value2 = 4999.3335
$ws.Range("K132").Value2 = 14998.0005
$ws.Range("M132").Value2 = -12468.0005

$ws = $wb.Worksheets.Item("LTW")
# Row 82 (anchor G82=12565) on LTW
$ws.Range("H82").Value2 = 407.2857
$ws.Range("I82").Value2 = 367
$ws.Range("J82").Value2 = 437.5
$ws.Range("K82").Value2 = 367
$ws.Range("L82").Value2 = 437.5
$ws.Range("M82").Value2 = -6
$ws.Range("N82").Value2 = -1159.5

# Row 85 (anchor G85=12565) on LTW
$ws.Range("H85").Value2 = 407.2857
$ws.Range("I85").Value2 = 367
$ws.Range("J85").Value2 = 437.5
$ws.Range("K85").Value2 = 367
$ws.Range("L85").Value2 = 437.5
$ws.Range("M85").Value2 = 881
$ws.Range("N85").Value2 = -2933.5

$ws = $wb.Worksheets.Item("WVR")
# Row 39 (anchor G39=3106) on WVR
$ws.Range("H39").Value2 = 0
$ws.Range("I39").Value2 = 0
$ws.Range("K39").Value2 = 0
$ws.Range("M39").Value2 = ""

# Row 93 (anchor G93=19613) on WVR
$ws.Range("H93").Value2 = 24000
$ws.Range("J93").Value2 = 24000
$ws.Range("L93").Value2 = 24000
$ws.Range("N93").Value2 = -28992

# Row 113 (anchor G113=27752) on WVR
$ws.Range("H113").Value2 = 582.84
$ws.Range("I113").Value2 = 1170.8572
$ws.Range("K113").Value2 = 3512.5716
$ws.Range("M113").Value2 = -1342.5716

# Row 132 (anchor G132=44029) on WVR
$ws.Range("H132").Value2 = 1295.7778
$ws.Range("I132").Value2 = 1295.7778
$ws.Range("K132").Value2 = 3887.3334
$ws.Range("M132").Value2 = -1357.3334
